$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.692.42'
$ws.Range('E2').Value = '  +1.63%  '
$ws.Range('D3').Value = '2.998.45'
$ws.Range('E3').Value = '  +3.56%  '
$ws.Range('D5').Value = '560.78'
$ws.Range('E5').Value = '  +2.05%  '
$ws.Range('D6').Value = '138.59'
$ws.Range('E6').Value = '  +13.34%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  +4.92%  '
$ws.Range('D9').Value = '2.987.33'
$ws.Range('E9').Value = '  +3.28%  '
$ws.Range('E10').Value = '  +5.65%  '
$ws.Range('E11').Value = '  +3.44%  '
$ws.Range('E12').Value = '  +4.72%  '
$ws.Range('E13').Value = '  +8.65%  '
$ws.Range('D14').Value = '33.63'
$ws.Range('E14').Value = '  +5.92%  '
$ws.Range('E15').Value = '  +3.26%  '
$ws.Range('D16').Value = '3.492.90'
$ws.Range('E16').Value = '  +3.65%  '
$ws.Range('D17').Value = '7.05'
$ws.Range('E17').Value = '  +8.93%  '
$ws.Range('D18').Value = '2.990.02'
$ws.Range('E18').Value = '  +3.55%  '
$ws.Range('D19').Value = '58.497.32'
$ws.Range('E19').Value = '  +1.46%  '
$ws.Range('D20').Value = '426.34'
$ws.Range('E20').Value = '  +4.54%  '
$ws.Range('D21').Value = '13.65'
$ws.Range('E21').Value = '  +5.93%  '
$ws.Range('D22').Value = '0.715'
$ws.Range('E22').Value = '  +9.07%  '
$ws.Range('D23').Value = '7.12'
$ws.Range('E23').Value = '  +5.42%  '
$ws.Range('D24').Value = '13.45'
$ws.Range('E24').Value = '  +7.03%  '
$ws.Range('D25').Value = '80.50'
$ws.Range('E25').Value = '  +4.53%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('E28').Value = '  +10.49%  '
$ws.Range('E29').Value = '  +3.44%  '
$ws.Range('D30').Value = '7.70'
$ws.Range('E30').Value = '  +6.74%  '
$ws.Range('D31').Value = '25.76'
$ws.Range('E31').Value = '  +4.62%  '
$ws.Range('D32').Value = '6.12'
$ws.Range('E32').Value = '  +1.53%  '
$ws.Range('D33').Value = '0.0986'
$ws.Range('E33').Value = '  +3.57%  '
$ws.Range('D34').Value = '5.78'
$ws.Range('E34').Value = '  +8.26%  '
$ws.Range('E35').Value = '  +6.02%  '
$ws.Range('D36').Value = '0.0₃0744'
$ws.Range('E36').Value = '  +20.09%  '
$ws.Range('D37').Value = '2.11'
$ws.Range('E37').Value = '  +4.21%  '
$ws.Range('D38').Value = '48.92'
$ws.Range('E38').Value = '  +1.30%  '
$ws.Range('D39').Value = '8.88'
$ws.Range('E39').Value = '  +4.92%  '
$ws.Range('D40').Value = '2.79'
$ws.Range('E40').Value = '  +17.57%  '
$ws.Range('D41').Value = '398.02'
$ws.Range('E41').Value = '  +11.33%  '
$ws.Range('D42').Value = '0.0352'
$ws.Range('E42').Value = '  +2.66%  '
$ws.Range('E43').Value = '  +2.84%  '
$ws.Range('D44').Value = '2.732.81'
$ws.Range('E44').Value = '  +4.79%  '
$ws.Range('D45').Value = '0.247'
$ws.Range('E45').Value = '  +8.12%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').Value = '125.51'
$ws.Range('E46').Value = '  +6.60%  '
$ws.Range('B47').Value = 'USDe'
$ws.Range('C47').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D47').Value = '0.999'
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('D48').Value = '2.03'
$ws.Range('E48').Value = '  +4.92%  '
$ws.Range('E49').Value = '  +2.89%  '
$ws.Range('D50').Value = '23.38'
$ws.Range('E50').Value = '  +3.36%  '
$ws.Range('B51').Value = 'Arweave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D51').Value = '31.28'
$ws.Range('E51').Value = '  +13.26%  '
